# BetaPhase 0.1.1 | Imagens e viewport prontos para animação
#
# This script reproduces, via Excel COM automation, the changes made to
# Calculadora.xlsx:
#   - H3 input changed from 2130 to 713 (formulas F3/G5/E7/G7 recalc)
#   - C7 input changed from 0.3 to 0.4
#   - The "label" cells (A5:B5, A7 "Em Porcentagem"/"Em Pixels", A9:B9, B7)
#     and the "value" cells (G5/H5, E7/F7/G7/H7) swap their visual styles
#     (font/fill/border/number format) with each other
#   - Active selection moves from A6:K6 to A8:K8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Stage the two style "looks" that need to be swapped using far-away
#    scratch cells, so the copy/paste doesn't clobber a source range
#    before every destination that needs it has been updated.
# ---------------------------------------------------------------------
$scratchLabelA = $ws.Range("Z1")   # look of A5 / B5  (border 4 side)
$scratchValueA = $ws.Range("Z2")   # look of G5 / E7 / G7 (border 4 side)
$scratchLabelB = $ws.Range("Z3")   # look of B7 (border 7 side)
$scratchValueB = $ws.Range("Z4")   # look of H5 / F7 / H7 (border 7 side)

$ws.Range("A5").Copy()
$scratchLabelA.PasteSpecial($xlPasteFormats)

$ws.Range("G5").Copy()
$scratchValueA.PasteSpecial($xlPasteFormats)

$ws.Range("B7").Copy()
$scratchLabelB.PasteSpecial($xlPasteFormats)

$ws.Range("H5").Copy()
$scratchValueB.PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Apply the "value"-look (previously used by G5/E7/G7) onto the
#    label cells A5, B5, A7, A9, B9; and the matching border-7 variant
#    onto... (none needed here, A side only for these)
# ---------------------------------------------------------------------
$scratchValueA.Copy()
$ws.Range("A5").PasteSpecial($xlPasteFormats)
$ws.Range("B5").PasteSpecial($xlPasteFormats)
$ws.Range("A7").PasteSpecial($xlPasteFormats)
$ws.Range("A9").PasteSpecial($xlPasteFormats)
$ws.Range("B9").PasteSpecial($xlPasteFormats)

$scratchValueB.Copy()
$ws.Range("B7").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 3) Apply the "label"-look (previously used by A5/B7) onto the value
#    cells G5, E7, G7 (border-4 variant) and H5, F7, H7 (border-7 variant)
# ---------------------------------------------------------------------
$scratchLabelA.Copy()
$ws.Range("G5").PasteSpecial($xlPasteFormats)
$ws.Range("E7").PasteSpecial($xlPasteFormats)
$ws.Range("G7").PasteSpecial($xlPasteFormats)

$scratchLabelB.Copy()
$ws.Range("H5").PasteSpecial($xlPasteFormats)
$ws.Range("F7").PasteSpecial($xlPasteFormats)
$ws.Range("H7").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) Clear the scratch area used for staging
# ---------------------------------------------------------------------
$ws.Range("Z1:Z4").Clear()

# ---------------------------------------------------------------------
# 5) Update the input values. Dependent formulas (F3, G5, E7, G7)
#    recalculate automatically.
# ---------------------------------------------------------------------
$ws.Range("H3").Value = 713
$ws.Range("C7").Value = 0.4

# ---------------------------------------------------------------------
# 6) Move the active selection from A6:K6 to A8:K8
# ---------------------------------------------------------------------
$ws.Range("A8:K8").Select()
